$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(7,  3,  "M",  12, "SL001-M",  1),
    @(8,  3,  "L",  12, "SL001-L",  1),
    @(9,  3,  "XL", 24, "SL001-XL", 1),
    @(10, 4,  "M",  12, "HD001-M",  1),
    @(11, 4,  "L",  24, "HD001-L",  1),
    @(12, 4,  "XL", 24, "HD001-XL", 1),
    @(13, 5,  "M",  24, "HD002-M",  1),
    @(14, 5,  "L",  24, "HD002-L",  1),
    @(15, 5,  "XL", 24, "HD002-XL", 1),
    @(16, 6,  "AS", 12, "BG001-AS", 1),
    @(17, 7,  "AS", 24, "BG002-AS", 1),
    @(18, 8,  "AS", 24, "SC001-AS", 1),
    @(19, 9,  "M",  12, "SR001-M",  1),
    @(20, 9,  "L",  12, "SR001-L",  1),
    @(21, 10, "L",  12, "JY001-L",  1)
)

$startRow = 8
$lastRow = $startRow + $data.Length - 1

# Fill column by column (A, B, C, D, E, F) to reproduce the authoring order
for ($col = 1; $col -le 6; $col++) {
    $r = $startRow
    foreach ($row in $data) {
        $ws.Cells.Item($r, $col).Value = $row[$col - 1]
        $r = $r + 1
    }
}

[void]$ws.Range("F14").Select()
